# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx price/volume/coin updates described in the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.185.62'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").Value = '1.804.65'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '223.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.553'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.41'
$ws.Range("D8").Style = "Normal"
$ws.Range("E9").Value = '  +2.23%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0721'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0928'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.87%  '
$ws.Range("D12").Value = '2.062.54'
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.811.45'
$ws.Range("E13").Value = '  +0.93%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.19%  '
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").Value = '34.213.74'
$ws.Range("E17").Value = '  -1.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.65'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '248.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.50%  '
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.97%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  -1.23%  '
$ws.Range("E24").Value = '  -0.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.37%  '
$ws.Range("E27").Value = '  +0.56%  '
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  +1.76%  '
$ws.Range("E31").Value = '  -0.49%  '
$ws.Range("E32").Value = '  +1.67%  '
$ws.Range("E33").Value = '  -1.11%  '
$ws.Range("E34").Value = '  -1.44%  '
$ws.Range("D35").Value = '1.419.20'
$ws.Range("E35").Value = '  -1.37%  '
$ws.Range("E36").Value = '  +2.74%  '
$ws.Range("E37").Value = '  +0.65%  '
$ws.Range("E38").Value = '  -1.58%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.946'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.47%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '80.69'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.40%  '
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.74'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.55%  '
$ws.Range("B42").Value = 'HuobiToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.36'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.63%  '
$ws.Range("E43").Value = '  +3.33%  '
$ws.Range("E44").Value = '  -1.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '108.18'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.83%  '
$ws.Range("E46").Value = '  +0.52%  '
$ws.Range("D47").Value = '1.962.57'
$ws.Range("E47").Value = '  +0.59%  '
$ws.Range("E48").Value = '  -1.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("E51").Value = '  +3.75%  '
